# cv125101a.xlsx — correção nos dados e inicio da analise PNAD 2009
#
# The original sheet had a two-row header where B2 carried a stray
# "unnamed: 1_level_1" label (an artifact of a pandas multi-index column
# flatten) and several section-title rows (sexo / cor ou raça / grupos de
# idade / nível de instrução / classes de rendimento mensal domiciliar per
# capita) that only held a label in column A with no data next to them.
# The footer rows (fonte/nota) were also dropped. This cleans all of that
# up: fix the header label and delete the now-unwanted rows so the data
# rows shift up into a contiguous block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the stray pandas header label in B2 -> "total"
$ws.Range("B2").Value = "total"

# Delete the section-header-only rows and the trailing source/footnote
# rows. Deleting from the bottom up keeps the remaining row numbers
# stable while we work.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

Write-Host "done"
